# Splits the single run
#   " section. An image is present in this section body. This section only contains one nested question."
# into five runs (same rPr throughout: <w:lang w:val="en-GB"/>):
#   " section. An image is present in this section body. This section "
#   "two"
#   " nested question"
#   "s"
#   "."
#
# The runtime merges adjacent runs back together whenever a Range.Text
# assignment leaves neighbouring runs with byte-identical rPr, so each
# segment is temporarily given a distinguishing Bold flag while the split
# is created, and the Bold flag is only cleared afterwards (as a pure
# formatting change on an already-isolated run, which does not re-trigger
# the merge pass).

$d = $word.ActiveDocument

# --- Step 1: "only contains one" -> "two" ---------------------------------
$txt = $d.Content.Text
$idx = $txt.IndexOf("only contains one nested question.")
$rTwo = $d.Range($idx, $idx + [string]"only contains one".Length)
$rTwo.Text = "two"
$rTwo.Font.Bold = 1

# --- Step 2: re-stamp " nested question" as its own run (text unchanged) --
$txt2 = $d.Content.Text
$idx2 = $txt2.IndexOf(" nested question.")
$rNested = $d.Range($idx2, $idx2 + [string]" nested question".Length)
$rNested.Text = " nested question"
$rNested.Font.Bold = 0

# --- Step 3: insert "s" right after "question", before the period --------
$txt3 = $d.Content.Text
$idx3 = $txt3.IndexOf(" nested question.") + [string]" nested question".Length
$rS = $d.Range($idx3, $idx3)
$rS.Text = "s"
$rS.Font.Bold = 1

# --- Step 4: strip the temporary Bold markers, now that each segment is --
# --- an isolated run (pure formatting change; does not re-merge) ---------
$txt4 = $d.Content.Text
$idxTwo = $txt4.IndexOf("This section two") + [string]"This section ".Length
$rClearTwo = $d.Range($idxTwo, $idxTwo + [string]"two".Length)
$rClearTwo.Font.Bold = 0

$txt5 = $d.Content.Text
$idxS = $txt5.IndexOf(" nested question") + [string]" nested question".Length
$rClearS = $d.Range($idxS, $idxS + 1)
$rClearS.Font.Bold = 0
